$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '41.539.77'
$ws.Range('E2').Value = '  +0.00%  '
$ws.Range('D3').Value = '2.495.82'
$ws.Range('E3').Value = '  +1.31%  '
$ws.Range('D4').Value = '0.995'
$ws.Range('E4').Value = '  -0.45%  '
$ws.Range('D5').Value = '311.63'
$ws.Range('E5').Value = '  +0.08%  '
$ws.Range('D6').Value = '93.10'
$ws.Range('E6').Value = '  -1.88%  '
$ws.Range('D7').Value = '0.540'
$ws.Range('E7').Value = '  -2.25%  '
$ws.Range('D8').Value = '0.998'
$ws.Range('E8').Value = '  -0.34%  '
$ws.Range('D9').Value = '0.495'
$ws.Range('E9').Value = '  -2.67%  '
$ws.Range('D10').Value = '32.39'
$ws.Range('E10').Value = '  -4.18%  '
$ws.Range('D11').Value = '0.0779'
$ws.Range('E11').Value = '  -0.44%  '
$ws.Range('E12').Value = '  +1.24%  '
$ws.Range('D13').Value = '2.867.30'
$ws.Range('E13').Value = '  +0.90%  '
$ws.Range('D14').Value = '6.83'
$ws.Range('E14').Value = '  -2.19%  '
$ws.Range('B15').Value = 'Chainlink'
$ws.Range('C15').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D15').Value = '15.36'
$ws.Range('E15').Value = '  +5.27%  '
$ws.Range('B16').Value = 'WrappedEther'
$ws.Range('C16').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D16').Value = '2.486.04'
$ws.Range('E16').Value = '  +0.79%  '
$ws.Range('D17').Value = '0.765'
$ws.Range('E17').Value = '  -3.07%  '
$ws.Range('D18').Value = '41.521.30'
$ws.Range('E18').Value = '  -0.04%  '
$ws.Range('D19').Value = '6.32'
$ws.Range('E19').Value = '  -0.97%  '
$ws.Range('D20').Value = '0.0₃0923'
$ws.Range('E20').Value = '  +0.69%  '
$ws.Range('D21').Value = '70.63'
$ws.Range('E21').Value = '  +1.66%  '
$ws.Range('D22').Value = '11.17'
$ws.Range('E22').Value = '  -3.78%  '
$ws.Range('D23').Value = '235.20'
$ws.Range('E23').Value = '  -0.76%  '
$ws.Range('D24').Value = '2.71'
$ws.Range('E24').Value = '  -2.51%  '
$ws.Range('E25').Value = '  -0.07%  '
$ws.Range('D26').Value = '1.90'
$ws.Range('E26').Value = '  -2.64%  '
$ws.Range('D27').Value = '24.52'
$ws.Range('E27').Value = '  -1.01%  '
$ws.Range('E28').Value = '  +1.04%  '
$ws.Range('D29').Value = '9.63'
$ws.Range('E29').Value = '  -1.22%  '
$ws.Range('D30').Value = '36.68'
$ws.Range('E30').Value = '  +0.91%  '
$ws.Range('D31').Value = '153.52'
$ws.Range('E31').Value = '  +0.21%  '
$ws.Range('D32').Value = '5.40'
$ws.Range('E32').Value = '  -4.36%  '
$ws.Range('D33').Value = '18.23'
$ws.Range('E33').Value = '  +6.23%  '
$ws.Range('E34').Value = '  -2.73%  '
$ws.Range('D35').Value = '0.0759'
$ws.Range('E35').Value = '  +0.38%  '
$ws.Range('D36').Value = '2.50'
$ws.Range('E36').Value = '  -1.96%  '
$ws.Range('D37').Value = '2.99'
$ws.Range('E37').Value = '  -0.91%  '
$ws.Range('D38').Value = '1.85'
$ws.Range('E38').Value = '  -1.93%  '
$ws.Range('D39').Value = '0.113'
$ws.Range('E39').Value = '  -1.36%  '
$ws.Range('D40').Value = '0.101'
$ws.Range('E40').Value = '  -3.68%  '
$ws.Range('D41').Value = '4.14'
$ws.Range('E41').Value = '  +2.52%  '
$ws.Range('D42').Value = '1.00'
$ws.Range('E42').Value = '  -0.40%  '
$ws.Range('D43').Value = '19.79'
$ws.Range('E43').Value = '  -7.81%  '
$ws.Range('D44').Value = '1.949.92'
$ws.Range('E44').Value = '  -1.87%  '
$ws.Range('D45').Value = '0.0283'
$ws.Range('E45').Value = '  -1.07%  '
$ws.Range('D46').Value = '2.97'
$ws.Range('E46').Value = '  -2.94%  '
$ws.Range('D47').Value = '8.73'
$ws.Range('E47').Value = '  +0.24%  '
$ws.Range('D48').Value = '2.722.96'
$ws.Range('E48').Value = '  +0.76%  '
$ws.Range('D49').Value = '95.99'
$ws.Range('E49').Value = '  -1.74%  '
$ws.Range('D50').Value = '0.176'
$ws.Range('E50').Value = '  -2.31%  '
$ws.Range('D51').Value = '67.15'
$ws.Range('E51').Value = '  -3.97%  '
